$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format first so numeric-looking strings (e.g. "578.64")
# are stored as text, matching the source inlineStr cells, not auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "67.556.51"
$ws.Range("E2").Value = "  +4.51%  "
$ws.Range("D3").Value = "3.260.28"
$ws.Range("E3").Value = "  +4.46%  "
$ws.Range("D5").Value = "578.64"
$ws.Range("E5").Value = "  +2.44%  "
$ws.Range("D6").Value = "182.29"
$ws.Range("E6").Value = "  +8.79%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").Value = "0.599"
$ws.Range("E8").Value = "  -1.28%  "
$ws.Range("D9").Value = "3.259.50"
$ws.Range("E9").Value = "  +4.21%  "
$ws.Range("E10").Value = "  +7.55%  "
$ws.Range("D11").Value = "6.72"
$ws.Range("E11").Value = "  +3.27%  "
$ws.Range("E12").Value = "  +7.23%  "
$ws.Range("D13").Value = "3.819.99"
$ws.Range("E13").Value = "  +3.92%  "
$ws.Range("E14").Value = "  +1.29%  "
$ws.Range("D15").Value = "28.54"
$ws.Range("E15").Value = "  +7.31%  "
$ws.Range("D16").Value = "67.528.63"
$ws.Range("E16").Value = "  +4.53%  "
$ws.Range("E17").Value = "  +4.84%  "
$ws.Range("D18").Value = "3.257.36"
$ws.Range("E18").Value = "  +3.75%  "
$ws.Range("E19").Value = "  +3.51%  "
$ws.Range("D20").Value = "13.58"
$ws.Range("E20").Value = "  +7.50%  "
$ws.Range("D21").Value = "375.77"
$ws.Range("E21").Value = "  +6.16%  "
$ws.Range("E22").Value = "  +6.37%  "
$ws.Range("E23").Value = "  +0.13%  "
$ws.Range("D24").Value = "71.23"
$ws.Range("E24").Value = "  +3.83%  "
$ws.Range("E25").Value = "  +4.42%  "
$ws.Range("E26").Value = "  +4.86%  "
$ws.Range("D27").Value = "9.58"
$ws.Range("E27").Value = "  +0.38%  "
$ws.Range("E28").Value = "  +2.79%  "
$ws.Range("E29").Value = "  +0.19%  "
$ws.Range("D30").Value = "5.72"
$ws.Range("E30").Value = "  +9.20%  "
$ws.Range("D31").Value = "1.97"
$ws.Range("E31").Value = "  +4.39%  "
$ws.Range("D32").Value = "22.72"
$ws.Range("E32").Value = "  +4.65%  "
$ws.Range("E33").Value = "  +0.02%  "
$ws.Range("E34").Value = "  +7.37%  "
$ws.Range("E35").Value = "  +6.34%  "
$ws.Range("D36").Value = "163.49"
$ws.Range("E36").Value = "  +3.65%  "
$ws.Range("E37").Value = "  +5.96%  "
$ws.Range("D38").Value = "0.852"
$ws.Range("E38").Value = "  +2.76%  "
$ws.Range("E39").Value = "  +6.05%  "
$ws.Range("D40").Value = "6.86"
$ws.Range("E40").Value = "  +13.35%  "
$ws.Range("D41").Value = "4.68"
$ws.Range("E41").Value = "  +13.76%  "
$ws.Range("E42").Value = "  +3.58%  "
$ws.Range("E43").Value = "  +8.90%  "
$ws.Range("D44").Value = "357.55"
$ws.Range("E44").Value = "  +13.61%  "
$ws.Range("D45").Value = "2.707.90"
$ws.Range("E45").Value = "  +2.44%  "
$ws.Range("D46").Value = "25.40"
$ws.Range("E46").Value = "  +7.33%  "
$ws.Range("D47").Value = "40.88"
$ws.Range("E47").Value = "  +4.12%  "
$ws.Range("D48").Value = "0.0681"
$ws.Range("E48").Value = "  +5.38%  "
$ws.Range("E49").Value = "  +4.39%  "
$ws.Range("D50").Value = "1.00"
$ws.Range("E50").Value = "  +8.32%  "
$ws.Range("E51").Value = "  +0.10%  "

# Restore the default "Normal" style on column D so no stray number-format style
# is left applied to the cells (keeps styles.xml equivalent to the original).
$ws.Range("D2:D51").Style = "Normal"
